$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (2014/12 IFRS column data refreshed)
$ws.Range("D2").Value = 1369
$ws.Range("E2").Value = 167
$ws.Range("F2").Value = 167
$ws.Range("G2").Value = 147
$ws.Range("H2").Value = 105
$ws.Range("I2").Value = 109
$ws.Range("J2").Value = -5
$ws.Range("K2").Value = 2457
$ws.Range("L2").Value = 899
$ws.Range("M2").Value = 1558
$ws.Range("N2").Value = 1561
$ws.Range("O2").Value = -3
$ws.Range("P2").Value = 65
$ws.Range("Q2").Value = 26
$ws.Range("R2").Value = 43
$ws.Range("S2").Value = -109
$ws.Range("T2").Value = 6
$ws.Range("U2").Value = 19
$ws.Range("V2").Value = 490
$ws.Range("W2").Value = 12.2
$ws.Range("X2").Value = 7.66
$ws.Range("Y2").Value = 7.1
$ws.Range("Z2").Value = 4.27
$ws.Range("AA2").Value = 57.73
$ws.Range("AB2").Value = 2526.91
$ws.Range("AC2").Value = 842
$ws.Range("AD2").Value = 9.52
$ws.Range("AE2").Value = 12756
$ws.Range("AF2").Value = 0.63
$ws.Range("AG2").Value = 130
$ws.Range("AH2").Value = 1.62
$ws.Range("AI2").Value = 14.53
$ws.Range("AJ2").Value = 12996741

# Row 3 (2015/12 IFRS column data refreshed)
$ws.Range("D3").Value = 1438
$ws.Range("E3").Value = 284
$ws.Range("F3").Value = 284
$ws.Range("G3").Value = 231
$ws.Range("H3").Value = 169
$ws.Range("I3").Value = 168
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 2352
$ws.Range("L3").Value = 630
$ws.Range("M3").Value = 1722
$ws.Range("N3").Value = 1725
$ws.Range("O3").Value = -3
$ws.Range("P3").Value = 65
$ws.Range("Q3").Value = 402
$ws.Range("R3").Value = -130
$ws.Range("S3").Value = -267
$ws.Range("T3").Value = 27
$ws.Range("U3").Value = 374
$ws.Range("V3").Value = 230
$ws.Range("W3").Value = 19.77
$ws.Range("X3").Value = 11.73
$ws.Range("Y3").Value = 10.24
$ws.Range("Z3").Value = 7.02
$ws.Range("AA3").Value = 36.57
$ws.Range("AB3").Value = 2764.84
$ws.Range("AC3").Value = 1294
$ws.Range("AD3").Value = 6.28
$ws.Range("AE3").Value = 14091
$ws.Range("AF3").Value = 0.58
$ws.Range("AG3").Value = 200
$ws.Range("AH3").Value = 2.46
$ws.Range("AI3").Value = 14.55
$ws.Range("AJ3").Value = 12996741

# Row 4 (2016/12 IFRS column data refreshed)
$ws.Range("D4").Value = 1439
$ws.Range("E4").Value = 355
$ws.Range("F4").Value = 355
$ws.Range("G4").Value = 331
$ws.Range("H4").Value = 256
$ws.Range("I4").Value = 256
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 2397
$ws.Range("L4").Value = 449
$ws.Range("M4").Value = 1948
$ws.Range("N4").Value = 1951
$ws.Range("O4").Value = -2
$ws.Range("P4").Value = 65
$ws.Range("Q4").Value = 512
$ws.Range("R4").Value = -209
$ws.Range("S4").Value = -224
$ws.Range("T4").Value = 5
$ws.Range("U4").Value = 507
$ws.Range("V4").Value = 30
$ws.Range("W4").Value = 24.65
$ws.Range("X4").Value = 17.81
$ws.Range("Y4").Value = 13.93
$ws.Range("Z4").Value = 10.79
$ws.Range("AA4").Value = 23.06
$ws.Range("AB4").Value = 3114.45
$ws.Range("AC4").Value = 1970
$ws.Range("AD4").Value = 8.07
$ws.Range("AE4").Value = 15935
$ws.Range("AF4").Value = 1
$ws.Range("AG4").Value = 330
$ws.Range("AH4").Value = 2.08
$ws.Range("AI4").Value = 15.78
$ws.Range("AJ4").Value = 12996741

# Row 5 (2017/12 IFRS column data refreshed)
$ws.Range("D5").Value = 1456
$ws.Range("E5").Value = 258
$ws.Range("F5").Value = 258
$ws.Range("G5").Value = 243
$ws.Range("H5").Value = 181
$ws.Range("I5").Value = 181
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 2563
$ws.Range("L5").Value = 472
$ws.Range("M5").Value = 2091
$ws.Range("N5").Value = 2092
$ws.Range("O5").Value = -1
$ws.Range("P5").Value = 65
$ws.Range("Q5").Value = 146
$ws.Range("R5").Value = -141
$ws.Range("S5").Value = -42
$ws.Range("T5").Value = 13
$ws.Range("U5").Value = 133
$ws.Range("V5").Value = 28
$ws.Range("W5").Value = 17.69
$ws.Range("X5").Value = 12.4
$ws.Range("Y5").Value = 8.949999999999999
$ws.Range("Z5").Value = 7.28
$ws.Range("AA5").Value = 22.58
$ws.Range("AB5").Value = 3330.53
$ws.Range("AC5").Value = 1391
$ws.Range("AD5").Value = 8.550000000000001
$ws.Range("AE5").Value = 17092
$ws.Range("AF5").Value = 0.7
$ws.Range("AG5").Value = 350
$ws.Range("AH5").Value = 2.94
$ws.Range("AI5").Value = 23.69
$ws.Range("AJ5").Value = 12996741

# Row 6 (2018/12 IFRS column data refreshed; J6/O6 intentionally absent, as before)
$ws.Range("D6").Value = 1619
$ws.Range("E6").Value = 103
$ws.Range("F6").Value = 103
$ws.Range("G6").Value = 60
$ws.Range("H6").Value = 23
$ws.Range("I6").Value = 24
$ws.Range("K6").Value = 2543
$ws.Range("L6").Value = 474
$ws.Range("M6").Value = 2069
$ws.Range("N6").Value = 2070
$ws.Range("P6").Value = 65
$ws.Range("Q6").Value = 175
$ws.Range("R6").Value = -302
$ws.Range("S6").Value = 2
$ws.Range("T6").Value = 76
$ws.Range("U6").Value = 99
$ws.Range("V6").Value = 76
$ws.Range("W6").Value = 6.38
$ws.Range("X6").Value = 1.44
$ws.Range("Y6").Value = 1.17
$ws.Range("Z6").Value = 0.91
$ws.Range("AA6").Value = 22.92
$ws.Range("AB6").Value = 3335.57
$ws.Range("AC6").Value = 188
$ws.Range("AD6").Value = 33.05
$ws.Range("AE6").Value = 16895
$ws.Range("AF6").Value = 0.37
$ws.Range("AG6").Value = 200
$ws.Range("AH6").Value = 3.22
$ws.Range("AI6").Value = 100.32
$ws.Range("AJ6").Value = 12996741

# Rows 7-9 (2019E/2020E/2021E): drop stale forecast figures, keep only rank/ticker/name
$ws.Range("D7:AJ7").ClearContents()
$ws.Range("D8:AJ8").ClearContents()
$ws.Range("D9:AJ9").ClearContents()
